$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 59.375
$ws.Range("I11").Value = 59.375
$ws.Range("K11").Value = 59.375
$ws.Range("M11").Value = 80.625
$ws.Range("H41").Value = 2719.2307
$ws.Range("I41").Value = 2434.5
$ws.Range("K41").Value = 2434.5
$ws.Range("M41").Value = -1994.5
$ws.Range("H53").Value = 3881.7856
$ws.Range("I53").Value = 1283
$ws.Range("K53").Value = 1283
$ws.Range("M53").Value = -646
$ws.Range("H107").Value = 5753.4
$ws.Range("I107").Value = 5611.4443
$ws.Range("J107").Value = 5966.3335
$ws.Range("K107").Value = 5611.4443
$ws.Range("L107").Value = 5966.3335
$ws.Range("M107").Value = -3691.4443
$ws.Range("N107").Value = -9806.333500000001
$ws.Range("H125").Value = 1275.5454
$ws.Range("I125").Value = 1010.6667
$ws.Range("J125").Value = 1374.875
$ws.Range("K125").Value = 9096.0003
$ws.Range("L125").Value = 12373.875
$ws.Range("M125").Value = -6636.0003
$ws.Range("N125").Value = -17293.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 22224616
$ws.Range("I74").Value = 27780270
$ws.Range("K74").Value = 27780270
$ws.Range("M74").Value = -27779396
$ws.Range("H77").Value = 22224616
$ws.Range("I77").Value = 27780270
$ws.Range("K77").Value = 138901350
$ws.Range("M77").Value = -138896982
$ws.Range("H96").Value = 41610.4
$ws.Range("J96").Value = 41610.4
$ws.Range("L96").Value = 41610.4
$ws.Range("N96").Value = -47102.4
$ws.Range("H97").Value = 1079.3043
$ws.Range("I97").Value = 1299
$ws.Range("K97").Value = 1299
$ws.Range("M97").Value = -803

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2749.3438
$ws.Range("I134").Value = 1105.8928
$ws.Range("J134").Value = 14253.5
$ws.Range("K134").Value = 3317.6784
$ws.Range("L134").Value = 42760.5
$ws.Range("M134").Value = -782.6784000000002
$ws.Range("N134").Value = -47830.5
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2402.6
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 4506.5
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 4506.5
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -5080.5
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10368
$ws.Range("H69").Value = 8826
$ws.Range("I69").Value = 8826
$ws.Range("K69").Value = 8826
$ws.Range("M69").Value = -8077
$ws.Range("H72").Value = 8826
$ws.Range("I72").Value = 8826
$ws.Range("K72").Value = 26478
$ws.Range("M72").Value = -22734
$ws.Range("H86").Value = 11092.714
$ws.Range("I86").Value = 4213.6665
$ws.Range("K86").Value = 4213.6665
$ws.Range("M86").Value = -3090.6665
$ws.Range("H89").Value = 11092.714
$ws.Range("I89").Value = 4213.6665
$ws.Range("K89").Value = 21068.3325
$ws.Range("M89").Value = -15452.3325
$ws.Range("H101").Value = 10000
$ws.Range("J101").Value = 10000
$ws.Range("L101").Value = 10000
$ws.Range("N101").Value = -16490
$ws.Range("H113").Value = 2402.6
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 4506.5
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 4506.5
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -8846.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 4000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = -12338
$ws.Range("H30").Value = 4000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 4000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = $null
$ws.Range("N30").Value = -12204
$ws.Range("H52").Value = 376.33334
$ws.Range("J52").Value = 376.33334
$ws.Range("L52").Value = 1129.00002
$ws.Range("N52").Value = -1661.00002
$ws.Range("H126").Value = 3749.5
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -2560
$ws.Range("N126").Value = -24877

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20073.188
$ws.Range("I93").Value = 17128
$ws.Range("J93").Value = 21054.916
$ws.Range("K93").Value = 17128
$ws.Range("L93").Value = 21054.916
$ws.Range("M93").Value = -15256
$ws.Range("N93").Value = -24798.916

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4173.6875
$ws.Range("I22").Value = 1641.1818
$ws.Range("J22").Value = 5500.2383
$ws.Range("K22").Value = 1641.1818
$ws.Range("L22").Value = 5500.2383
$ws.Range("M22").Value = -1346.1818
$ws.Range("N22").Value = -6090.2383
$ws.Range("H27").Value = 4173.6875
$ws.Range("I27").Value = 1641.1818
$ws.Range("J27").Value = 5500.2383
$ws.Range("K27").Value = 1641.1818
$ws.Range("L27").Value = 5500.2383
$ws.Range("M27").Value = -1534.1818
$ws.Range("N27").Value = -5714.2383
$ws.Range("H100").Value = 11289.556
$ws.Range("I100").Value = 8682.923000000001
$ws.Range("J100").Value = 13710
$ws.Range("K100").Value = 8682.923000000001
$ws.Range("L100").Value = 13710
$ws.Range("M100").Value = -8141.923000000001
$ws.Range("N100").Value = -14792

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 97613
$ws.Range("J16").Value = 97613
$ws.Range("L16").Value = 97613
$ws.Range("N16").Value = -98197
$ws.Range("H120").Value = 77183.5
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").Value = $null
$ws.Range("H123").Value = 52000
$ws.Range("J123").Value = 52000
$ws.Range("L123").Value = 52000
$ws.Range("N123").Value = -61800
$ws.Range("H126").Value = 3561.2144
$ws.Range("I126").Value = 3408.375
$ws.Range("J126").Value = 3765
$ws.Range("K126").Value = 10225.125
$ws.Range("L126").Value = 11295
$ws.Range("M126").Value = -7755.125
$ws.Range("N126").Value = -16235
$ws.Range("H131").Value = 88750
$ws.Range("J131").Value = 88750
$ws.Range("L131").Value = 88750
$ws.Range("N131").Value = -98830
$ws.Range("H141").Value = 69313.875
$ws.Range("J141").Value = 69313.875
$ws.Range("L141").Value = 69313.875
$ws.Range("N141").Value = -79673.875
